$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81 - JUANITO A. ANGCAYA
$ws.Range("B81").Value = 45093
$ws.Range("C81").Value = "his"
$ws.Range("D81").Value = "MR"
$ws.Range("F81").Value = "JUANITO"
$ws.Range("G81").Value = "ANGCAYA"
$ws.Range("E81").Value = "ANGCAYA"
$ws.Range("I81").Value = "Admin Aide I"
$ws.Range("J81").Value = "Picnic Grove"
$ws.Range("L81").Value = 45062
$ws.Range("O81").Value = 77.485999999999947
$ws.Range("P81").Value = 340.65
$ws.Range("N81").Value = 12092
$ws.Range("T81").Value = "his Compulsory retirement"

# Row 82 - ADELAIDA C. LUCIANO
$ws.Range("B82").Value = 45093
$ws.Range("C82").Value = "her"
$ws.Range("D82").Value = "MS"
$ws.Range("F82").Value = "ADELAIDA"
$ws.Range("G82").Value = "CREUS"
$ws.Range("E82").Value = "LUCIANO"
$ws.Range("I82").Value = "Community Affairs Asst II"
$ws.Range("J82").Value = "Community/Brgy. Affair's Office"
$ws.Range("L82").Value = 45088
$ws.Range("O82").Value = 58.746000000000038
$ws.Range("P82").Value = 267.98200000000003
$ws.Range("N82").Value = 18212
$ws.Range("R82").Value = 0.048192699999999998
$ws.Range("N82").Value = 18212
$ws.Range("T82").Value = "his Compulsory retirement"

$ws.Range("K81").Select
